# "Generate Report for Archive"
# The localization status report is regenerated: the "zh-cn"/"de-de" status
# on the Overview sheet (and the matching per-language sheets) moves from
# "Ready for handoff" to "In Translation", and the now-narrower Status
# column is resized to fit the shorter text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"
$narrowWidth = 12.5   # character-width input that yields the narrower "Status" column seen in the regenerated report

# --- Overview sheet: Status shows in both the zh-cn (col E) and de-de (col F) columns ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $narrowWidth
$wsOverview.Columns.Item(6).ColumnWidth = $narrowWidth

# --- Per-language detail sheets: Status is column C ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $narrowWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $narrowWidth
